$d = $word.ActiveDocument

# wdFindContinue = 1
$wdFindContinue = 1
# wdColorRed = 255 (RGB(255,0,0))
$wdColorRed = 255

# --- 1) ". 가급적 리턴함수를 ... 필요하다." paragraph:
#    color "가급적 리턴함수를 사용하여 명확히 구현하는 것이 필요" red,
#    leaving ". " and "하다." uncolored ---
$r1 = $d.Content
$r1.Find.Execute("가급적 리턴함수를 사용하여 명확히 구현하는 것이 필요", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$r1.Font.Color = $wdColorRed

# --- 2) "f4(5) 이건 가능할까? -> f4(p2=5) 로 호출" paragraph:
#    color "f4(5)" red ---
$r2 = $d.Content
$r2.Find.Execute("f4(5)", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$r2.Font.Color = $wdColorRed

#    color "p2=5" (inside the second "f4(p2=5)") red ---
$r2b = $d.Content
$r2b.Find.Execute("p2=5", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$r2b.Font.Color = $wdColorRed

# --- 3) "f5<- function(...) { print(...) ... }" paragraph:
#    color only the "..." argument list red, keep "function(" and ") " uncolored ---
$r3 = $d.Content
$r3.Find.Execute("function(...)", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$dotsRange = $d.Range($r3.Start + 9, $r3.Start + 12)
$dotsRange.Font.Color = $wdColorRed
